$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rename the "Durée" (duration) coded values from bare "5" / "10" to the
# more explicit "5mn" / "10mn" across every sheet that references them
# (Echant1_impl, Echant2_impl, Params1_expl, Params2_expl).
# ---------------------------------------------------------------------------

# Echant1_impl
$ws1 = $wb.Worksheets.Item("Echant1_impl")
$ws1.Range("D2").Value = "5mn"
$ws1.Range("D3").Value = "10mn"

# Echant2_impl
$ws2 = $wb.Worksheets.Item("Echant2_impl")
$ws2.Range("D2").Value = "5mn"
$ws2.Range("D3").Value = "10mn"

# Params1_expl
$ws4 = $wb.Worksheets.Item("Params1_expl")
foreach ($r in 2..6)   { $ws4.Cells.Item($r, 4).Value = "5mn" }
foreach ($r in 7..16)  { $ws4.Cells.Item($r, 4).Value = "10mn" }
foreach ($r in 17..19) { $ws4.Cells.Item($r, 4).Value = "5mn" }
foreach ($r in 20..22) { $ws4.Cells.Item($r, 4).Value = "10mn" }

# Params2_expl
$ws5 = $wb.Worksheets.Item("Params2_expl")
foreach ($r in 2..5)   { $ws5.Cells.Item($r, 4).Value = "5mn" }
foreach ($r in 6..9)   { $ws5.Cells.Item($r, 4).Value = "10mn" }
foreach ($r in 10..12) { $ws5.Cells.Item($r, 4).Value = "5mn" }

# ---------------------------------------------------------------------------
# Cosmetic: update the remembered selection on each sheet.
# Echant1_impl (the tab-selected sheet) is reselected last so it remains the
# active sheet/tab after the script finishes.
# ---------------------------------------------------------------------------

$ws2.Range("D2:D3").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("Modl_impl")
$ws3.Range("F29").Select() | Out-Null

$ws4.Range("F24").Select() | Out-Null

$ws5.Range("D7:D9").Select() | Out-Null

$ws1.Activate()
$ws1.Range("F18").Select() | Out-Null
